$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new research-hours row (row 20), matching the date-formatted style of A19
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A20").Value = 46079
$ws.Range("B20").Value = "worked on the hypothesis - U=M/C"
$ws.Range("C20").Value = 1

# Match the selection state saved with the workbook
$ws.Range("C20").Select()
